$d = $word.ActiveDocument

# 1. Update the title text
$d.Content.Find.Execute("Apex, Testing And Debugging", $false, $false, $false, $false, $false, $true, 1, $false, "Apex Specialist SuperBadge", 2)
